$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 3
    4  = 4
    5  = 1
    6  = 3
    7  = 2
    8  = 0
    9  = 5
    10 = 3
    11 = 3
    12 = 3
    13 = 1
    14 = 0
    15 = 3
    16 = 2
    17 = 3
    18 = 6
    19 = 2
    20 = 1
    21 = 2
    22 = 6
    23 = 2
    24 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
